$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 1.68
$ws.Range("G2").Value = 2
$ws.Range("I2").Value = 6.6
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1.34
$ws.Range("Q2").Value = 1.94
$ws.Range("V2").Value = 1.18
$ws.Range("W2").Value = 2
$ws.Range("G3").Value = 1.39
$ws.Range("H3").Value = 10
$ws.Range("P3").Value = 2.04
$ws.Range("R3").Value = 1.4
$ws.Range("U3").Value = 1.71
$ws.Range("W3").Value = 3.45
$ws.Range("Z3").Value = 110
$ws.Range("AA3").Value = 460
$ws.Range("AI3").Value = 170
$ws.Range("AO3").Value = 310
$ws.Range("G4").Value = 4.4
$ws.Range("I4").Value = 2.24
$ws.Range("J4").Value = 3.3
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 2.94
$ws.Range("O4").Value = 1.45
$ws.Range("P4").Value = 1.65
$ws.Range("Q4").Value = 2.34
$ws.Range("R4").Value = 1.24
$ws.Range("S4").Value = 4.6
$ws.Range("T4").Value = 1.98
$ws.Range("U4").Value = 1.88
$ws.Range("Y4").Value = 7.8
$ws.Range("AB4").Value = 15
$ws.Range("F5").Value = 2.12
$ws.Range("G5").Value = 2.24
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 4.1
$ws.Range("K5").Value = 3.65
$ws.Range("M5").Value = 1.09
$ws.Range("N5").Value = 3.15
$ws.Range("P5").Value = 1.76
$ws.Range("Q5").Value = 2.14
$ws.Range("R5").Value = 1.28
$ws.Range("T5").Value = 1.9
$ws.Range("U5").Value = 1.98
$ws.Range("W5").Value = 1.68
$ws.Range("X5").Value = 14.5
$ws.Range("Y5").Value = 13.5
$ws.Range("AB5").Value = 8.6
$ws.Range("AD5").Value = 17
$ws.Range("AE5").Value = 55
$ws.Range("AO5").Value = 65
$ws.Range("F6").Value = 2.4
$ws.Range("G6").Value = 2.58
$ws.Range("H6").Value = 3.3
$ws.Range("I6").Value = 3.65
$ws.Range("J6").Value = 3.15
$ws.Range("K6").Value = 3.25
$ws.Range("N6").Value = 2.78
$ws.Range("P6").Value = 1.6
$ws.Range("R6").Value = 1.22
$ws.Range("V6").Value = 1.37
$ws.Range("W6").Value = 1.63
$ws.Range("Z6").Value = 24
$ws.Range("AA6").Value = 75
$ws.Range("AB6").Value = 9.4
$ws.Range("AD6").Value = 16
$ws.Range("AF6").Value = 15
$ws.Range("AG6").Value = 12.5
$ws.Range("AJ6").Value = 42
$ws.Range("AK6").Value = 34
$ws.Range("F7").Value = 2.42
$ws.Range("G7").Value = 2.64
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 3.65
$ws.Range("V7").Value = 1.38
$ws.Range("W7").Value = 1.61
$ws.Range("X7").Value = 12
$ws.Range("Y7").Value = 11
$ws.Range("Z7").Value = 24
$ws.Range("AA7").Value = 70
$ws.Range("AB7").Value = 8.4
$ws.Range("AC7").Value = 7.4
$ws.Range("AD7").Value = 16
$ws.Range("AE7").Value = 980
$ws.Range("AF7").Value = 980
$ws.Range("AG7").Value = 12.5
$ws.Range("AH7").Value = 980
$ws.Range("AI7").Value = 70
$ws.Range("AJ7").Value = 980
$ws.Range("AK7").Value = 980
$ws.Range("AL7").Value = 60
$ws.Range("AM7").Value = 180
$ws.Range("S8").Value = 5.5
$ws.Range("W8").Value = 1.59
$ws.Range("AM8").Value = 230
$ws.Range("F9").Value = 2.28
$ws.Range("G9").Value = 2.44
$ws.Range("J9").Value = 3.45
$ws.Range("M9").Value = 1.09
$ws.Range("N9").Value = 3.2
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 1.76
$ws.Range("Q9").Value = 2.14
$ws.Range("R9").Value = 1.28
$ws.Range("S9").Value = 4
$ws.Range("T9").Value = 1.89
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.38
$ws.Range("W9").Value = 1.7
$ws.Range("X9").Value = 980
$ws.Range("Y9").Value = 980
$ws.Range("Z9").Value = 980
$ws.Range("AA9").Value = 70
$ws.Range("AB9").Value = 9
$ws.Range("AC9").Value = 7.8
$ws.Range("AD9").Value = 980
$ws.Range("AF9").Value = 980
$ws.Range("AG9").Value = 980
$ws.Range("AH9").Value = 980
$ws.Range("AI9").Value = 60
$ws.Range("AM9").Value = 140
$ws.Range("G10").Value = 3.15
$ws.Range("H10").Value = 2.7
$ws.Range("I10").Value = 2.8
$ws.Range("P10").Value = 1.58
$ws.Range("Q10").Value = 2.5
$ws.Range("U10").Value = 1.83
$ws.Range("W10").Value = 1.46
$ws.Range("AL10").Value = 75
$ws.Range("AN10").Value = 65
$ws.Range("AO10").Value = 1000
$ws.Range("F11").Value = 2.12
$ws.Range("G11").Value = 2.38
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 3.85
$ws.Range("J11").Value = 3.45
$ws.Range("M11").Value = 1.06
$ws.Range("O11").Value = 1.29
$ws.Range("P11").Value = 1.94
$ws.Range("Q11").Value = 1.88
$ws.Range("T11").Value = 1.72
$ws.Range("V11").Value = 1.36
$ws.Range("W11").Value = 1.72
$ws.Range("Y11").Value = 15
$ws.Range("Z11").Value = 26
$ws.Range("AA11").Value = 960
$ws.Range("AB11").Value = 11
$ws.Range("AC11").Value = 9.800000000000001
$ws.Range("AD11").Value = 16.5
$ws.Range("AE11").Value = 44
$ws.Range("AG11").Value = 13.5
$ws.Range("AH11").Value = 21
$ws.Range("AI11").Value = 55
$ws.Range("AK11").Value = 34
$ws.Range("AN11").Value = 21
$ws.Range("AO11").Value = 46
